# Daily attendance processing - 2025-10-22 07:43:39
# Applies the recorded-by reordering, summary count corrections, column widen,
# and "Pending" -> "Not Recorded" restyle to the Session Analysis Results sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

# --- Widen the "Students" helper column (column I / index 9) ---
# (the engine's ColumnWidth<->stored-width mapping carries a constant
#  +5/6 padding, so back it out to land on an exact stored width of 14)
$ws.Columns.Item(9).ColumnWidth = 14 - (5/6)

# --- Reorder "Recorded By" values: move the leading "System" entry to the end ---
$recordedByRows = 2,4,5,7,8,29,31,32,34,35,56,58,59,61,62,83,84,85,109,110,111,135,136,137
foreach ($r in $recordedByRows) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2
    if ($val) {
        $parts = $val -split ", "
        $reordered = (@($parts[1..($parts.Length - 1)]) + @($parts[0])) -join ", "
        $cell.Value = $reordered
    }
}

# --- Updated Missing/Pending session counters on the summary block ---
$ws.Cells.Item(7, 12).Value = 3
$ws.Cells.Item(8, 12).Value = 39

# --- Updated per-group coverage counters (P/Q columns) for rows 18-20 ---
foreach ($r in 18,19,20) {
    $ws.Cells.Item($r, 16).Value = 1
    $ws.Cells.Item($r, 17).Value = 7
}

# --- Rows 101, 127, 153: still-pending sessions are now flagged "Not Recorded"
#     with the pink highlight instead of the yellow "Pending" highlight ---
$notRecordedRows = 101,127,153
foreach ($r in $notRecordedRows) {
    $rowRange = $ws.Range($ws.Cells.Item($r,1), $ws.Cells.Item($r,9))
    $rowRange.Interior.Color = 12695295
    $rowRange.Font.Color = 0
    $rowRange.HorizontalAlignment = -4108
    $rowRange.VerticalAlignment = -4108
    $ws.Cells.Item($r, 9).Value = "Not Recorded"
}
